# community.xlsx - "test case added" commit
# Fills in the Procedure / Test Condition / Expected Result columns for the
# new "Web Sustainability" community test-case rows (rows 3-11) and updates the
# active-sheet selection to match the author's final cursor position (B24).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3: "User views community detail page." already has Test Condition/
# Expected Result; only the Procedure (D3) cell is new. ---
$ws.Range("D3").Value = '1-Open http://localhost:8080/index.html
2-Login with username emre.gurer@gmail.com and password 123.
3-Click on details of Web Sustainability Community'

# --- Row 4: "User views details of a topic." ---
$ws.Range("B4").Value = 'User views details of a topic.'
$ws.Range("D4").Value = '1- Click on "Sustainable Web Design Resource" topic.'
$ws.Range("C4").Value = 'Topic details with comments and resources are listed.'

# --- Row 5: "User returns to community detail page." ---
$ws.Range("B5").Value = 'User returns to community detail page.'
$ws.Range("D5").Value = '1- Click on return to community on topic detail page.'
$ws.Range("C5").Value = '"Web Sustainability" community''s detail page is opened.'

# --- Row 6: "User views member profile." ---
$ws.Range("B6").Value = 'User views member profile.'
$ws.Range("D6").Value = '1- Click on İlyas Alper Karatepe text'
$ws.Range("C6").Value = 'User profile page is opened.'

# --- Row 7: "User views all members" ---
$ws.Range("B7").Value = 'User views all members'
$ws.Range("D7").Value = '1- Click on see all members'
$ws.Range("C7").Value = 'List is displayed:
İlyas Alper Karatepe
Emre Gürer'

# --- Row 8: "User views upcoming event in detail" ---
$ws.Range("B8").Value = 'User views upcoming event in detail'
$ws.Range("D8").Value = '1- Click on meeting on 21 December'
$ws.Range("C8").Value = 'Meeting page is opened.'

# --- Row 9: "User views resources of community" ---
$ws.Range("B9").Value = 'User views resources of community'
$ws.Range("D9").Value = '1- Click on alper.png'
$ws.Range("C9").Value = 'Image is opened'

# --- Row 10: "User views community requests" ---
$ws.Range("B10").Value = 'User views community requests'
$ws.Range("D10").Value = '1- Click on community requests'
$ws.Range("C10").Value = 'No request is displayed.'

# --- Row 11: "User creates a topic" ---
$ws.Range("B11").Value = 'User creates a topic'
$ws.Range("D11").Value = '1- Click on create topic
2-Title about web
3-Description web will be explained
4-tag web
5-Click create'
$ws.Range("C11").Value = 'Topic is created and topic detail page is opened.'

# --- Row heights grew to fit the newly wrapped multi-line text. ---
$ws.Rows.Item(3).RowHeight = 75
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(11).RowHeight = 75

# --- Author left the selection on B24 before saving. ---
$ws.Range("B24").Select()

